$d = $word.ActiveDocument

$replacements = @(
    @("479÷3=159, 2", "490÷8=61, 2"),
    @("167÷6=27, 5", "768÷7=109, 5"),
    @("698÷4=174, 2", "220÷8=27, 4"),
    @("208÷2=104, 0", "189÷9=21, 0"),
    @("995÷7=142, 1", "610÷2=305, 0"),
    @("429÷5=85, 4", "350÷7=50, 0"),
    @("148÷8=18, 4", "181÷5=36, 1"),
    @("894÷3=298, 0", "533÷7=76, 1"),
    @("365÷9=40, 5", "316÷2=158, 0"),
    @("516÷2=258, 0", "633÷4=158, 1"),
    @("305÷4=76, 1", "190÷6=31, 4"),
    @("611÷5=122, 1", "501÷3=167, 0"),
    @("682÷8=85, 2", "120÷5=24, 0"),
    @("790÷4=197, 2", "925÷7=132, 1"),
    @("394÷3=131, 1", "914÷6=152, 2"),
    @("285÷8=35, 5", "404÷8=50, 4"),
    @("548÷9=60, 8", "333÷8=41, 5"),
    @("530÷9=58, 8", "562÷8=70, 2"),
    @("356÷2=178, 0", "418÷3=139, 1"),
    @("862÷8=107, 6", "876÷4=219, 0"),
    @("199÷8=24, 7", "373÷7=53, 2"),
    @("134÷4=33, 2", "586÷7=83, 5"),
    @("837÷7=119, 4", "564÷8=70, 4"),
    @("949÷8=118, 5", "512÷3=170, 2"),
    @("101÷7=14, 3", "300÷9=33, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
